# "Changed customer under test"
# The workbook's per-sheet "Customer" cell is updated from
# "TrustwaveEP Automation #3" to "TrustwaveEP Automation #4" on the sheets
# that carry that field, and the remembered cell-selection on several sheets
# is moved to reflect where the user was working.

$wb = $excel.ActiveWorkbook

$newCustomer = "TrustwaveEP Automation #4"

# --- DownloadFromPortalTest: Customer is column B (row 2) ---
$wsDownload = $wb.Worksheets.Item("DownloadFromPortalTest")
$wsDownload.Range("B2").Value = $newCustomer
$wsDownload.Activate()
$wsDownload.Range("B2").Select()

# --- VerifyEndPointOkAtPortalTest: Customer is column B (row 2) ---
$wsVerify = $wb.Worksheets.Item("VerifyEndPointOkAtPortalTest")
$wsVerify.Range("B2").Value = $newCustomer
$wsVerify.Activate()
$wsVerify.Range("B2").Select()

# --- ClientLogToPortalTest: Customer is column C (row 2); this sheet stays
#     the active/selected tab, so it is handled last so its activation and
#     selection remain the ones persisted for the workbook. ---
$wsClientLog = $wb.Worksheets.Item("ClientLogToPortalTest")
$wsClientLog.Range("C2").Value = $newCustomer
$wsClientLog.Activate()
$wsClientLog.Range("C2:C3").Select()

# Best-effort: reflect the updated Excel window geometry recorded for the
# workbook view.
$win = $excel.ActiveWindow
$win.Left = -108
$win.Top = -108
$win.Width = 23256
$win.Height = 12576
